{"js": "// Rename the report title from \"R\u00e9union - 16/10/14\" to \"R\u00e9union 1 - 16/10/14\"\n// (this report becomes \"report 1\" now that a second report is being added,\n// per the commit message \"Ajout du rapport 2\").\nconst body = context.document.body;\nconst results = body.search(\"R\u00e9union - 16/10/14\", { matchCase: true, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertText(\"R\u00e9union 1 - 16/10/14\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Rename the report title from \"R\u00e9union - 16/10/14\" to \"R\u00e9union 1 - 16/10/14\"\n# (this report becomes \"report 1\" now that a second report is being added,\n# per the commit message \"Ajout du rapport 2\").\n$d = $word.ActiveDocument\n$find = $d.Content.Find\n$find.Text = \"R\u00e9union - 16/10/14\"\n$find.Replacement.Text = \"R\u00e9union 1 - 16/10/14\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2)\n"}
